$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; insert it as a new row
# right before the existing row 101, shifting rows 101-147 down to
# 102-148 (matches the diff: dimension grows from A1:T147 to A1:T148,
# and every data row from 101 on is the prior row's values).
$ws.Rows.Item(101).EntireRow.Insert()

$ws.Cells.Item(101, 1).Value  = 11
$ws.Cells.Item(101, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(101, 3).Value  = "Bíobío"
$ws.Cells.Item(101, 4).Value  = 44917
$ws.Cells.Item(101, 5).Value  = 8
$ws.Cells.Item(101, 6).Value  = "Fruta"
$ws.Cells.Item(101, 7).Value  = 100108
$ws.Cells.Item(101, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value  = 100108002
$ws.Cells.Item(101, 10).Value = "Mango"
$ws.Cells.Item(101, 11).Value = "Sin especificar"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 200
$ws.Cells.Item(101, 14).Value = 7000
$ws.Cells.Item(101, 15).Value = 7500
$ws.Cells.Item(101, 16).Value = 7250
$ws.Cells.Item(101, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(101, 18).Value = "Ecuador"
$ws.Cells.Item(101, 19).Value = 1812
$ws.Cells.Item(101, 20).Value = 4
